$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 100.07
$ws.Range("H2").Value = 0.4
$ws.Range("I2").Value = 100.07
$ws.Range("J2").Value = 0.33
$ws.Range("K2").Value = 18.94
$ws.Range("L2").Value = 0.4
$ws.Range("M2").Value = 19.28

$ws.Range("E3").Value = 100.07
$ws.Range("H3").Value = 0.63
$ws.Range("I3").Value = 100.07
$ws.Range("J3").Value = 0.52
$ws.Range("K3").Value = 11.05
$ws.Range("L3").Value = 0.63
$ws.Range("M3").Value = 11.63

$ws.Range("E4").Value = 100.07
$ws.Range("I4").Value = 100.07
$ws.Range("K4").Value = 15.38
$ws.Range("M4").Value = 15.8

$ws.Range("E5").Value = 100.07
$ws.Range("I5").Value = 100.07
$ws.Range("K5").Value = 16.72
$ws.Range("M5").Value = 17.1

$ws.Range("E6").Value = 100.07
$ws.Range("H6").Value = 0.49
$ws.Range("I6").Value = 100.07
$ws.Range("J6").Value = 0.41
$ws.Range("K6").Value = 17.14
$ws.Range("L6").Value = 0.49
$ws.Range("M6").Value = 17.52

$ws.Range("E7").Value = 100.07
$ws.Range("I7").Value = 100.07
$ws.Range("K7").Value = 11.33
$ws.Range("M7").Value = 11.89

$ws.Range("E8").Value = 100.07
$ws.Range("I8").Value = 100.07
$ws.Range("K8").Value = 17.04
$ws.Range("M8").Value = 17.41

$ws.Range("E9").Value = 100.07
$ws.Range("I9").Value = 100.07
$ws.Range("K9").Value = 15.8
$ws.Range("M9").Value = 16.21

$ws.Range("E10").Value = 100.07
$ws.Range("H10").Value = 0.65
$ws.Range("I10").Value = 100.07
$ws.Range("J10").Value = 0.54
$ws.Range("K10").Value = 12.37
$ws.Range("L10").Value = 0.65
$ws.Range("M10").Value = 12.89

$ws.Range("E11").Value = 100.07
$ws.Range("I11").Value = 100.07
$ws.Range("K11").Value = 16.59
$ws.Range("M11").Value = 16.98

$ws.Range("E12").Value = 100.07
$ws.Range("I12").Value = 100.07
$ws.Range("K12").Value = 19.11
$ws.Range("M12").Value = 19.45

$ws.Range("E13").Value = 100.07
$ws.Range("H13").Value = 0.27
$ws.Range("I13").Value = 100.07
$ws.Range("J13").Value = 0.22
$ws.Range("K13").Value = 30.95
$ws.Range("L13").Value = 0.27
$ws.Range("M13").Value = 31.16

$ws.Range("E14").Value = 100.07
$ws.Range("H14").Value = 0.59
$ws.Range("I14").Value = 100.07
$ws.Range("J14").Value = 0.49
$ws.Range("K14").Value = 13.61
$ws.Range("L14").Value = 0.59
$ws.Range("M14").Value = 14.08
